# Lesson 2.2 Dividing into cases -- apply author's edits
#
# 1) Slide 1 (title slide): subtitle's "Lesson 1.7" -> "Lesson 2.2"
# 2) Slide 8 ("Another example"): insert two new comment lines
#    (";; GIVEN: ..." and ";; RETURNS: ...") right after the
#    ";; ball-after-tick : Ball -> Ball" line, before the
#    ";; STRATEGY: ..." lines.

$p = $ppt.ActivePresentation

# --- 1) Slide 1: update the lesson number in the subtitle -----------------
$slide1 = $p.Slides.Item(1)
$subtitle = $slide1.Shapes.Item(2)
$subtitleRange = $subtitle.TextFrame.TextRange

for ($i = 1; $i -le $subtitleRange.Paragraphs().Count; $i++) {
    $para = $subtitleRange.Paragraphs($i)
    if ($para.Text -match "Lesson") {
        $para.Text = "Lesson 2.2"
    }
}

# --- 2) Slide 8: add GIVEN/RETURNS lines to the code comment block --------
$slide8 = $p.Slides.Item(8)
$content8 = $slide8.Shapes.Item(2)
$contentRange = $content8.TextFrame.TextRange

$firstPara = $contentRange.Paragraphs(1)
$firstPara.InsertAfter("`r;; GIVEN: The state of a ball b`r;; RETURNS: the state of given ball at the next tick")
